$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 header text changes from "Temperature" to "temperature" (lowercase).
# C1 ("HomeOdds") and D1 ("OverUnder") keep their text.
$ws.Range("B1").Value = "temperature"

# Update data rows: column B becomes 60, columns C and D are removed (odds data dropped).
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 2).Value = 60
    $ws.Cells.Item($r, 3).ClearContents()
    $ws.Cells.Item($r, 4).ClearContents()
}

# Column A width (best fit) set explicitly to match the recorded width.
$ws.Columns.Item(1).ColumnWidth = 17.5

# Update selection to C2
$ws.Range("C2").Select()

# Adjust window size to match target view state
$excel.ActiveWindow.Width = 8280
$excel.ActiveWindow.Height = 19120
